# Weekly update: a new price observation is inserted as row 32, pushing the
# existing rows (formerly 32-114) down by one (now 33-115).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(32).Insert()

$ws.Cells.Item(32,1).Value2 = 3
$ws.Cells.Item(32,2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(32,3).Value2 = "Coquimbo"
$ws.Cells.Item(32,4).Value2 = 44414
$ws.Cells.Item(32,5).Value2 = 5
$ws.Cells.Item(32,6).Value2 = 100112001
$ws.Cells.Item(32,7).Value2 = "Berenjena"
$ws.Cells.Item(32,8).Value2 = "Sin especificar"
$ws.Cells.Item(32,9).Value2 = "Primera"
$ws.Cells.Item(32,10).Value2 = 50
$ws.Cells.Item(32,11).Value2 = 12000
$ws.Cells.Item(32,12).Value2 = 12000
$ws.Cells.Item(32,13).Value2 = 12000
$ws.Cells.Item(32,14).Value2 = "`$/caja 60 unidades"
$ws.Cells.Item(32,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(32,16).Value2 = 200
$ws.Cells.Item(32,17).Value2 = 60
$ws.Cells.Item(32,18).Value2 = "Hortaliza"
